$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: Educators -> Education
$ws.Range("B2").Value = "Education"

# B15: Port safety -> Port Safety
$ws.Range("B15").Value = "Port Safety"

# Row 16: was "Other" / "Another group nit listed" -> becomes "Not Specified" row (shifted from old row 17)
$ws.Range("A16").Value = 16
$ws.Range("B16").Value = "Not Specified"
$ws.Range("C16").Value = "Comment does not include reference to specifc user or uses"
$ws.Range("D16").Value = 19

# Row 17: becomes "Resilience Planning" row (shifted from old row 18)
$ws.Range("A17").Value = 17
$ws.Range("B17").Value = "Resilience Planning"
$ws.Range("C17").Value = "Resilience planning, identification of risks, and risk reduction"
$ws.Range("D17").Value = 14

# Row 18: new "Shoreline Protection" row
$ws.Range("A18").Value = 18
$ws.Range("B18").Value = "Shoreline Protection"
$ws.Range("C18").Value = "Erosion protection, living shorelines, etc."
$ws.Range("D18").Value = 11
